$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "26.866.55"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -1.97%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.806.46"
$cell.Style = "Normal"

$ws.Range("E4").Value = "  +0.08%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "310.23"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.21%  "

$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("E7").Value = "  +2.80%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3740"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.39%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07363"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.85%  "

$ws.Range("E10").Value = "  -1.39%  "

$ws.Range("E11").Value = "  -3.60%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "5.348"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -1.18%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "6.528"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.44%  "

$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.07050"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "1.715.93"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -5.90%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "91.03"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -3.19%  "

$ws.Range("E17").Value = "  +0.14%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000008734"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("E19").Value = "  +0.11%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "14.74"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -2.85%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "26.874.16"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -1.94%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.308"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("E23").Value = "  -1.40%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.968.42"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -4.28%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "1.911"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -2.39%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "151.28"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "18.45"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.70%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.151"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -8.96%  "

$ws.Range("E29").Value = "  -1.57%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "115.85"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -1.89%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.08900"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.57%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.7703"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.06%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.158"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -3.62%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.475"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -1.27%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.905"
$cell.Style = "Normal"

$ws.Range("E36").Value = "  +0.08%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.115"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +0.26%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01958"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -2.01%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.05240"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.58%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.414"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +4.73%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "7.250"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -1.66%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.5344"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.33%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.901"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.60%  "

$ws.Range("E44").Value = "  -3.68%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "8.550"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.14%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.5059"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.81%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "10.29"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -3.78%  "

$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.10%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "103.88"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -1.78%  "

$ws.Range("E50").Value = "  -2.33%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.06316"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -1.11%  "
